$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (far outside the used range) used to force text-typed values
# for new strings that would otherwise be auto-converted to numbers by Excel
# when assigned directly via .Value (mirrors typing a Text-formatted value and
# pasting just the value, leaving no extra style on the destination cell).
$helper = $ws.Cells.Item(200, 200)
$helper.NumberFormat = "@"

function Set-TextValue {
    param($Target, [string]$Text)
    $helper.Value = $Text
    $helper.Copy()
    $Target.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

$ws.Range("D2").Value = '49.332.17'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").Value = '2.622.75'
$ws.Range("E3").Value = '  -0.87%  '
Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  +0.01%  '
Set-TextValue $ws.Range("D5") '111.58'
$ws.Range("E5").Value = '  -2.21%  '
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("E7").Value = '  -1.15%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -1.15%  '
Set-TextValue $ws.Range("D10") '39.36'
$ws.Range("E10").Value = '  -4.27%  '
Set-TextValue $ws.Range("D11") '19.97'
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("E12").Value = '  -1.30%  '
$ws.Range("E13").Value = '  +1.33%  '
Set-TextValue $ws.Range("D14") '7.53'
$ws.Range("E14").Value = '  +2.25%  '
$ws.Range("D15").Value = '3.038.00'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '2.631.47'
$ws.Range("E16").Value = '  -0.66%  '
Set-TextValue $ws.Range("D17") '0.851'
$ws.Range("E17").Value = '  -2.03%  '
$ws.Range("D18").Value = '49.298.01'
$ws.Range("E18").Value = '  -0.99%  '
Set-TextValue $ws.Range("D19") '13.23'
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("E20").Value = '  -1.72%  '
Set-TextValue $ws.Range("D21") '6.65'
$ws.Range("E21").Value = '  -1.96%  '
$ws.Range("E22").Value = '  -1.16%  '
Set-TextValue $ws.Range("D23") '267.72'
$ws.Range("E23").Value = '  -3.48%  '
Set-TextValue $ws.Range("D24") '68.94'
$ws.Range("E24").Value = '  -4.33%  '
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D26") '0.999'
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D27") '25.94'
$ws.Range("E27").Value = '  -3.14%  '
$ws.Range("E28").Value = '  +1.57%  '
Set-TextValue $ws.Range("D29") '2.19'
$ws.Range("E29").Value = '  -1.45%  '
Set-TextValue $ws.Range("D30") '0.138'
$ws.Range("E30").Value = '  -2.09%  '
Set-TextValue $ws.Range("D31") '34.50'
$ws.Range("E31").Value = '  -4.54%  '
Set-TextValue $ws.Range("D32") '49.50'
$ws.Range("E32").Value = '  -1.67%  '
$ws.Range("E33").Value = '  +0.61%  '
Set-TextValue $ws.Range("D34") '0.0805'
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("E35").Value = '  -0.16%  '
Set-TextValue $ws.Range("D36") '18.96'
$ws.Range("E36").Value = '  -3.07%  '
Set-TextValue $ws.Range("D37") '4.95'
$ws.Range("E37").Value = '  +1.50%  '
$ws.Range("E38").Value = '  -2.46%  '
$ws.Range("E39").Value = '  -0.02%  '
Set-TextValue $ws.Range("D40") '129.07'
$ws.Range("E40").Value = '  +2.33%  '
$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D41") '2.29'
$ws.Range("E41").Value = '  +2.08%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D42") '22.73'
$ws.Range("E42").Value = '  +3.31%  '
$ws.Range("E43").Value = '  -1.07%  '
$ws.Range("E44").Value = '  +5.80%  '
$ws.Range("D45").Value = '2.057.60'
$ws.Range("E45").Value = '  -1.07%  '
Set-TextValue $ws.Range("D46") '3.26'
$ws.Range("E46").Value = '  -1.78%  '
Set-TextValue $ws.Range("D47") '2.13'
$ws.Range("E47").Value = '  +7.97%  '
$ws.Range("E48").Value = '  -5.65%  '
$ws.Range("E49").Value = '  -2.91%  '
$ws.Range("E50").Value = '  -3.58%  '
Set-TextValue $ws.Range("D51") '58.38'
$ws.Range("E51").Value = '  -3.30%  '

# Clean up helper cell and clipboard marching-ants state
$helper.Clear()
$excel.CutCopyMode = $false

